$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.944.27'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.383.41'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.14'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.98'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.61'
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('E10').Value = '  -1.04%  '
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '3.963.59'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.79'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '3.392.51'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = '61.056.28'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('E18').Value = '  -3.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.64'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.97'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '382.94'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.88'
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000116'
$ws.Range('E25').Value = '  -2.98%  '
$ws.Range('D26').Value = '3.523.96'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.26'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.99'
$ws.Range('E30').Value = '  -1.72%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -4.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.24'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.95'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.73'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').Value = '3.416.10'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.99'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('E39').Value = '  -4.10%  '
$ws.Range('E40').Value = '  -1.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.92'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D47').Value = '2.458.86'
$ws.Range('E47').Value = '  -3.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.00'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.71'
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('E50').Value = '  +8.74%  '
$ws.Range('E51').Value = '  +1.54%  '
